$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计" (i.e. before the
#    sheet that is currently in position 2, "2022-Q3").
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-fetch a fresh reference to the "2022-Q3" sheet (sheet references
# captured before a worksheet-collection change can resolve to the wrong
# sheet afterwards, so grab it again now that the collection is settled).
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Copy header formatting (bold / border / centered - style used by the other
# quarterly sheets) from "2022-Q3" row 1 into the new sheet, then set the
# header text.
$q3Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the index-column formatting (style used for column A in the other
# quarterly sheets) for the two data rows.
$q3Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Row 2: 诺德量化蓝筹增强混合C
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'005083"
$newSheet.Range("C2").Value = "诺德量化蓝筹增强混合C"
$newSheet.Range("D2").Value = "'0.54"
$newSheet.Range("E2").Value = "'86.59"
$newSheet.Range("F2").Value = "'4.24"
$newSheet.Range("G2").Value = "'0.0229"
$newSheet.Range("H2").Value = 2

# Row 3: 诺德量化蓝筹增强混合A
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'005082"
$newSheet.Range("C3").Value = "诺德量化蓝筹增强混合A"
$newSheet.Range("D3").Value = "'0.00"
$newSheet.Range("E3").Value = "'86.59"
$newSheet.Range("F3").Value = "'4.24"
$newSheet.Range("G3").Value = 0
$newSheet.Range("H3").Value = 2

# ---------------------------------------------------------------------------
# 2) Add a "2022-Q4" row to the "总计" summary sheet, right after the header
#    row, pushing the existing data rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Copy the formatting of the row that just got pushed down to row 3 into the
# new row 2 (keeps column A's centered/bold style and leaves B:D unstyled,
# matching the rest of the sheet).
$totalSheet.Range("A3:D3").Copy()
$totalSheet.Range("A2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.02

# Fix up the running index in column A for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
